# RMA Complete Flow (Issue Credit) - SO TO RMA Receipt To Create Credit Memo
# Update the "RMA Details Maintenance Grid" sheet with a new set of RMA
# numbers / Salesforce record ids (rows 2-4: RMA header name, RMA line
# name, and the associated Salesforce Id), replacing the previous
# "RMA-PU7Q-*" test data with a fresh "RMA-Q190-*" batch.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RMA Details Maintenance Grid")

# Row 2
$ws.Range("E2").Value = "RMA-Q190-001"
$ws.Range("F2").Value = "RMA-Q190-1-1"
$ws.Range("J2").Value = "a7s5f000000xNY8AAM"

# Row 3
$ws.Range("E3").Value = "RMA-Q190-002"
$ws.Range("F3").Value = "RMA-Q190-1-2"
$ws.Range("J3").Value = "a7s5f000000xNY9AAM"

# Row 4
$ws.Range("E4").Value = "RMA-Q190-003"
$ws.Range("F4").Value = "RMA-Q190-1-3"
$ws.Range("J4").Value = "a7s5f000000xNYAAA2"
